$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay plain text (matches the source file's
    # inline-string cells) even when the new value looks numeric,
    # e.g. "1.00" or "0.480" -- Excel would otherwise coerce those
    # into numbers and drop the significant trailing zeros.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Update the Price (D) and Volume(1h) (E) columns row by row ---

# Row 2
Set-TextValue $ws.Range("D2") '61.417.45'
Set-TextValue $ws.Range("E2") '  -4.40%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.318.79'
Set-TextValue $ws.Range("E3") '  -4.93%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.05%  '

# Row 5
Set-TextValue $ws.Range("D5") '567.56'
Set-TextValue $ws.Range("E5") '  -3.48%  '

# Row 6
Set-TextValue $ws.Range("D6") '128.87'
Set-TextValue $ws.Range("E6") '  -3.02%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +0.03%  '

# Row 8
Set-TextValue $ws.Range("D8") '3.318.03'
Set-TextValue $ws.Range("E8") '  -4.90%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.480'
Set-TextValue $ws.Range("E9") '  -0.68%  '

# Row 10
Set-TextValue $ws.Range("D10") '7.30'
Set-TextValue $ws.Range("E10") '  -4.99%  '

# Row 11
Set-TextValue $ws.Range("E11") '  -3.77%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.378'
Set-TextValue $ws.Range("E12") '  -1.97%  '

# Row 13
Set-TextValue $ws.Range("D13") '3.890.30'
Set-TextValue $ws.Range("E13") '  -4.95%  '

# Row 14
Set-TextValue $ws.Range("E14") '  -0.23%  '

# Row 15
Set-TextValue $ws.Range("D15") '3.325.85'
Set-TextValue $ws.Range("E15") '  -4.78%  '

# Row 16
Set-TextValue $ws.Range("E16") '  -5.14%  '

# Row 17
Set-TextValue $ws.Range("D17") '24.78'
Set-TextValue $ws.Range("E17") '  +1.79%  '

# Row 18
Set-TextValue $ws.Range("D18") '61.619.04'
Set-TextValue $ws.Range("E18") '  -4.06%  '

# Row 19
Set-TextValue $ws.Range("D19") '13.61'
Set-TextValue $ws.Range("E19") '  +1.28%  '

# Row 20
Set-TextValue $ws.Range("D20") '5.70'
Set-TextValue $ws.Range("E20") '  -0.61%  '

# Row 21
Set-TextValue $ws.Range("D21") '9.03'
Set-TextValue $ws.Range("E21") '  -9.23%  '

# Row 22
Set-TextValue $ws.Range("D22") '354.82'
Set-TextValue $ws.Range("E22") '  -7.68%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.559'
Set-TextValue $ws.Range("E23") '  -3.14%  '

# Row 24
Set-TextValue $ws.Range("D24") '0.998'
Set-TextValue $ws.Range("E24") '  -0.23%  '

# Row 25
Set-TextValue $ws.Range("D25") '3.451.74'
Set-TextValue $ws.Range("E25") '  -4.95%  '

# Row 26
Set-TextValue $ws.Range("D26") '69.64'
Set-TextValue $ws.Range("E26") '  -6.39%  '

# Row 27
Set-TextValue $ws.Range("D27") '0.0000108'
Set-TextValue $ws.Range("E27") '  -5.39%  '

# Row 28
Set-TextValue $ws.Range("D28") '1.00'
Set-TextValue $ws.Range("E28") '  +0.26%  '

# Row 29
Set-TextValue $ws.Range("E29") '  -0.03%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.45'
Set-TextValue $ws.Range("E30") '  -1.61%  '

# Row 31
Set-TextValue $ws.Range("D31") '7.88'
Set-TextValue $ws.Range("E31") '  -2.11%  '

# Row 32
Set-TextValue $ws.Range("E32") '  -5.49%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.150'
Set-TextValue $ws.Range("E33") '  -2.32%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -0.01%  '

# Row 35
Set-TextValue $ws.Range("D35") '3.347.25'
Set-TextValue $ws.Range("E35") '  -4.89%  '

# Row 36
Set-TextValue $ws.Range("D36") '22.66'
Set-TextValue $ws.Range("E36") '  -1.88%  '

# Row 37
Set-TextValue $ws.Range("D37") '5.34'
Set-TextValue $ws.Range("E37") '  -0.40%  '

# Row 38
Set-TextValue $ws.Range("D38") '6.84'
Set-TextValue $ws.Range("E38") '  +0.01%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.49'
Set-TextValue $ws.Range("E39") '  -2.35%  '

# Row 40
Set-TextValue $ws.Range("D40") '161.42'
Set-TextValue $ws.Range("E40") '  -1.56%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.0763'
Set-TextValue $ws.Range("E41") '  -2.16%  '

# Row 42
Set-TextValue $ws.Range("D42") '1.00'
Set-TextValue $ws.Range("E42") '  +0.16%  '

# Row 43
Set-TextValue $ws.Range("D43") '4.36'
Set-TextValue $ws.Range("E43") '  -0.06%  '

# Row 44
Set-TextValue $ws.Range("D44") '41.04'
Set-TextValue $ws.Range("E44") '  -1.74%  '

# Row 45
Set-TextValue $ws.Range("E45") '  -7.55%  '

# Row 46
Set-TextValue $ws.Range("E46") '  -4.75%  '

# Row 47
Set-TextValue $ws.Range("E47") '  -4.56%  '

# Row 50
Set-TextValue $ws.Range("D50") '0.863'
Set-TextValue $ws.Range("E50") '  -6.10%  '

# Row 51
Set-TextValue $ws.Range("D51") '21.43'
Set-TextValue $ws.Range("E51") '  +3.39%  '

# --- Rows 48/49 swap places: Cosmos now ranks above EnergySwap ---
# (Rank numbers in column A are unaffected; only Coin/Link/Price/Volume move.)
Set-TextValue $ws.Range("B48") "Cosmos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D48") "6.74"
Set-TextValue $ws.Range("E48") "  -0.34%  "

Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "22.09"
Set-TextValue $ws.Range("E49") "  -8.15%  "
